$d = $word.ActiveDocument

function Set-BulletParagraphText($doc, [string]$oldText, [string]$newText) {
    # Replace a "What we like/don't like" bullet paragraph's text while
    # preserving its exact run layout (a leading empty run followed by a
    # text run), matching the source document's authoring pattern.
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $text = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($text -eq $oldText) {
            $snippet = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'
            $para.Range.InsertXML($snippet) | Out-Null
            return $true
        }
    }
    return $false
}

# 1. Update the title (appears both as Heading1 at the top and as a bold
#    paragraph near the bottom) - replace every occurrence.
$d.Content.Find.Execute(
    "Play Dragon & Phoenix Slot for Free - Betsoft 2019 Game", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Play Dragon & Phoenix Slot Game - Free Review", 2) | Out-Null

# 2. "What we like" bullet list: insert a brand-new bullet
#    "Asian-themed graphics and theme" right before "Regular and expanding
#    wilds", matching the existing run layout (an empty run followed by a
#    text run) by inserting raw OOXML for the new paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Regular and expanding wilds") {
        $para.Range.InsertParagraphBefore() | Out-Null
        $insertedRange = $d.Paragraphs.Item($i).Range
        $snippet = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Asian-themed graphics and theme</w:t></w:r></w:p>'
        $insertedRange.InsertXML($snippet) | Out-Null
        break
    }
}

# 3. "Free spins with increasing multipliers" -> "Scatters and free spins feature"
Set-BulletParagraphText $d "Free spins with increasing multipliers" "Scatters and free spins feature" | Out-Null

# 4. Remove the old "Asian-themed graphics and symbols" bullet entirely
#    (its content moved to the new bullet added in step 2) and shorten the
#    "Potential for big wins..." bullet.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Asian-themed graphics and symbols") {
        $para.Range.Delete()
        break
    }
}

Set-BulletParagraphText $d "Potential for big wins up to 8098x total bet" "Potential for big wins" | Out-Null

# 5. "What we don't like" bullet list: swap / rewrite the two bullets.
Set-BulletParagraphText $d "RTP is only 96.01%" "Medium-high volatility" | Out-Null
Set-BulletParagraphText $d "Medium-high volatility may not suit all players" "RTP of 96.01%" | Out-Null

# 6. Update the meta-description style italic paragraph near the end.
$d.Content.Find.Execute(
    "Read our Dragon and Phoenix slot review and play the game for free. Features, pros, and cons, and potential for big wins up to 8098x total bet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Dragon & Phoenix slot game for free and discover its Asian-themed graphics, wilds, scatters, and potential for big wins.",
    2) | Out-Null
